$d = $word.ActiveDocument

foreach ($p in @($d.Paragraphs)) {
    $t = $p.Range.Text
    if ($t -like "After years of working for two software development companies*") {
        $p.Range.Delete()
        break
    }
}
